$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original "alcohol" measurement sheet has an extra column (N) that
# duplicates/derives from column M. Remove column M entirely, which shifts
# the old column N left into the M position.
[void]$ws.Columns.Item(13).Delete()

# Excel leaves the active cell on the column that now occupies the deleted
# column's position after a column delete.
[void]$ws.Range("M1").Select()
